$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.440382361412048
$ws.Range("B1").Value = 5.685306549072266
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.692667245864868
$ws.Range("E1").Value = 1.608485221862793
